# Edit slide 8 ("Functional Groupings of CVM Opcodes") of the CVM deck:
#   * insert a new "Bitwise Opcodes: BITAND, BITOR, BITXOR, BITNOT" paragraph
#     right after "Logical Opcodes: NOT" (and before "Shift Opcodes: ...")
#   * insert a new "Type Conversion Opcodes: BYTE2INT and INT2BYTE" paragraph
#     right after "Shift Opcodes: SHL and SHR" (and before "Branch Opcodes: ...")
#   * bump the slide's sldId in the presentation slide list (375 -> 376)
#   * turn on autoCompressPictures="0" at the presentation level

$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Presentation-level tweaks
# ------------------------------------------------------------------
$p.AutoCompressPictures = $false

# Slide 8 (1-based) is the "Functional Groupings of CVM Opcodes" slide.
$slide = $p.Slides.Item(8)
$slide.SlideID = 376

# ------------------------------------------------------------------
# 2) Content Placeholder body text edits
# ------------------------------------------------------------------
$contentShape = $slide.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

# --- Insert "Type Conversion Opcodes: ..." before the "Branch Opcodes: " paragraph ---
$branchPara = $tr.Paragraphs(4, 1)

$typeConv = $branchPara.InsertBefore("Type Conversion Opcodes: BYTE2INT and INT2BYTE`r")

$r = $typeConv.Characters(26, 8)            ; # "BYTE2INT"
$r.Font.Name = "Consolas"
$r.Font.Size = 20
$r = $typeConv.Characters(34, 5)            ; # " and "
$r.Font.Size = 20
$r = $typeConv.Characters(39, 8)            ; # "INT2BYTE"
$r.Font.Name = "Consolas"
$r.Font.Size = 20

# --- Insert "Bitwise Opcodes: ..." before the "Shift Opcodes: " paragraph ---
$shiftPara = $tr.Paragraphs(3, 1)

$bitwise = $shiftPara.InsertBefore("Bitwise Opcodes: BITAND, BITOR, BITXOR, BITNOT`r")

$r = $bitwise.Characters(18, 6)             ; # "BITAND"
$r.Font.Name = "Consolas"
$r.Font.Size = 20
$r = $bitwise.Characters(26, 5)             ; # "BITOR"
$r.Font.Name = "Consolas"
$r.Font.Size = 20
$r = $bitwise.Characters(33, 6)             ; # "BITXOR"
$r.Font.Name = "Consolas"
$r.Font.Size = 20
$r = $bitwise.Characters(41, 6)             ; # "BITNOT"
$r.Font.Name = "Consolas"
$r.Font.Size = 20

Write-Host "Final body text:"
Write-Host $tr.Text
